# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the freshly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2,5,7,9,10,27,35
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 203
$ws1.Range("F5").Value = 18816
$ws1.Range("F7").Value = 282
$ws1.Range("F9").Value = 7000
$ws1.Range("F10").Value = 441
$ws1.Range("F27").Value = 10
$ws1.Range("F35").Value = 12264

# Sheet "全部类型" (all types) - rows 2,5,7,9,10,27,37
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 203
$ws4.Range("F5").Value = 18816
$ws4.Range("F7").Value = 282
$ws4.Range("F9").Value = 7000
$ws4.Range("F10").Value = 441
$ws4.Range("F27").Value = 10
$ws4.Range("F37").Value = 12264

$wb.Save()
